# Weapon editor and selector updated
# Append new "Advantages" rows to Sheet1 / Table1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$penalty20 = "O personagem, quando em combate corpo-a-corpo, consegue usar seus oponentes comocobertura, bloqueando e atrapalhando os ataques de inimigos. Sempre que estiver enfrentando mais de um oponente ao mesmo tempo, o personagem é capaz de se posicionar de modo a atrapalhar os ataques dos outros inimigos.`nCada um dos oponentes recebe uma penalidade de 20% em ataque. Não pode ser usado contra armas de`nlonga distância."
$penalty30 = "O personagem, quando em combate corpo-a-corpo, consegue usar seus oponentes comocobertura, bloqueando e atrapalhando os ataques de inimigos. Sempre que estiver enfrentando mais de um oponente ao mesmo tempo, o personagem é capaz de se posicionar de modo a atrapalhar os ataques dos outros inimigos.`nCada um dos oponentes recebe uma penalidade de 30% em ataque. Não pode ser usado contra armas de`nlonga distância."
$amorVerdadeiro = "Você ama alguém do fundo do seu coração, e não importa o que ou quem tente interferir neste relacionamento, nada abalará o amor mútuo entre vocês sentem um pelo outro. Independente da situação, você vai ter forças para defender seu amor.`nEm toda a situação que a pessoa que você ama estiver em apuros, o personagem recebe um bônus especial de +10% em todos os seus Testes de Atributos Físicos e Perícias que forem realizados na tentativa de ajudá-la. também torna todos os seus Testes de WILL Fáceis contra Sedução."
$aparenciaInofensiva = "Você não aparenta ser perigoso. Na verdade, os oponentes menosprezam sua presença, ninguém acredita que você seja capaz de realizar algum grande feito e dificilmente o tomarão como uma ameaça (a não ser aqueles que já conhecem o Personagem há tempo suficiente para saber se isso é verdade ou não). Considere que o personagem automaticamente vence a Iniciativa quando um combate começar.`nPorém, depois de enfrentar um adversário uma vez, ele poderá já estar ciente do seu potencial e`nnão lhe permitirá nenhuma vantagem (então a Iniciativa deve ser testada normalmente)."

# Shared-string order matters: names for the two "Alvo Elusivo" rows were
# entered first, then their descriptions were pasted back in reverse
# (30% before 20%), before continuing with the remaining rows in order.
$ws.Cells.Item(72, 1).Value = "Alvo Elusivo (1)"
$ws.Cells.Item(73, 1).Value = "Alvo Elusivo (2)"
$ws.Cells.Item(73, 4).Value = $penalty30
$ws.Cells.Item(72, 4).Value = $penalty20

$ws.Cells.Item(72, 2).Value = 1
$ws.Cells.Item(73, 2).Value = 2

$ws.Cells.Item(74, 1).Value = "Amor verdadeiro"
$ws.Cells.Item(74, 2).Value = 1
$ws.Cells.Item(74, 4).Value = $amorVerdadeiro

$ws.Cells.Item(75, 1).Value = "Aparência Inofensiva"
$ws.Cells.Item(75, 2).Value = 1
$ws.Cells.Item(75, 4).Value = $aparenciaInofensiva

$ws.Cells.Item(76, 1).Value = "Ataque"

# Match the row heights used by the author for the wrapped description text.
$ws.Rows.Item(72).RowHeight = 90
$ws.Rows.Item(73).RowHeight = 90
$ws.Rows.Item(74).RowHeight = 90
$ws.Rows.Item(75).RowHeight = 90

# Grow the table / list object to include the new rows.
$wb.Worksheets.Item("Sheet1").ListObjects.Item("Table1").Resize($ws.Range("A1:D76"))

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("D4").Select()
